$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold/border/centered) from H1 onto the two new
# header cells before setting their text, so I1/J1 match the existing
# header formatting (style index 1 in the original workbook).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns (I = I0, J = IF) for each data row
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 5

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 5

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 6

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 5

$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 4

$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 6

$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 5
